$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.898
$ws.Range("C14").Value = -12.592
$ws.Range("C21").Value = -12.45
$ws.Range("D22").Value = -7.617
$ws.Range("C23").Value = -12.766
$ws.Range("D24").Value = -7.693
$ws.Range("C25").Value = -11.642
$ws.Range("C26").Value = -11.912
$ws.Range("D28").Value = -8.175000000000001
$ws.Range("C29").Value = -12.129
$ws.Range("D36").Value = -7.417999999999999
$ws.Range("D45").Value = -7.582000000000001
$ws.Range("D48").Value = -7.457000000000001
$ws.Range("D49").Value = -8.062999999999999
$ws.Range("D52").Value = -7.99
$ws.Range("C53").Value = -11.611
$ws.Range("D53").Value = -7.567
$ws.Range("D54").Value = -8.077
$ws.Range("C57").Value = -13.742
$ws.Range("C59").Value = -12.336
$ws.Range("C69").Value = -10.921
$ws.Range("D70").Value = -7.419
$ws.Range("C79").Value = -12.471
$ws.Range("C83").Value = -13.351
$ws.Range("D86").Value = -8.265999999999998
$ws.Range("D87").Value = -7.927
$ws.Range("D89").Value = -7.869
$ws.Range("C91").Value = -12.72
$ws.Range("C93").Value = -10.836
$ws.Range("D101").Value = -7.556999999999999
$ws.Range("C103").Value = -12.088
